$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row labels: "_old" -> "_FV2304", "_new" -> "_FV2310" ---
$headers = @(
    "Segmentname_FV2304","Segmentgruppe_FV2304","Segment_FV2304","Datenelement_FV2304",
    "Segment ID_FV2304","Code_FV2304","Qualifier_FV2304","Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304","Bedingung_FV2304","diff",
    "Segmentname_FV2310","Segmentgruppe_FV2310","Segment_FV2310","Datenelement_FV2310",
    "Segment ID_FV2310","Code_FV2310","Qualifier_FV2310","Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310","Bedingung_FV2310"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value2 = $headers[$i]
}

# --- Turn the data range into an Excel Table (ListObject) ---
$tableRange = $ws.Range("A1:U58")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# --- Freeze the header row (split at row 2 / pane frozen) ---
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
